$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for the Price column so numeric-looking strings
# (e.g. "0.9985", "29.578.23") are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.578.23"
$ws.Range("E2").Value = "  -0.73%  "
$ws.Range("D3").Value = "1.856.07"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("D4").Value = "0.9985"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "243.77"
$ws.Range("E5").Value = "  -0.59%  "
$ws.Range("D6").Value = "0.6438"
$ws.Range("E6").Value = "  +0.45%  "
$ws.Range("D7").Value = "0.9993"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "0.3018"
$ws.Range("E8").Value = "  +1.20%  "
$ws.Range("D9").Value = "0.07522"
$ws.Range("E9").Value = "  +0.30%  "
$ws.Range("D10").Value = "24.41"
$ws.Range("E10").Value = "  +1.06%  "
$ws.Range("D11").Value = "0.07663"
$ws.Range("E11").Value = "  -0.25%  "
$ws.Range("D12").Value = "1.915.86"
$ws.Range("E12").Value = "  +2.81%  "
$ws.Range("D13").Value = "5.050"
$ws.Range("D14").Value = "0.6904"
$ws.Range("E14").Value = "  +0.77%  "
$ws.Range("D15").Value = "83.94"
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("D16").Value = "0.000009581"
$ws.Range("E16").Value = "  +1.18%  "
$ws.Range("D17").Value = "6.227"
$ws.Range("E17").Value = "  +2.74%  "
$ws.Range("D18").Value = "2.158.48"
$ws.Range("E18").Value = "  +2.42%  "
$ws.Range("D19").Value = "29.600.07"
$ws.Range("E19").Value = "  -0.55%  "
$ws.Range("D20").Value = "237.93"
$ws.Range("E20").Value = "  -0.71%  "
$ws.Range("D21").Value = "12.62"
$ws.Range("E21").Value = "  -0.60%  "
$ws.Range("D24").Value = "0.9998"
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("D25").Value = "157.34"
$ws.Range("E25").Value = "  -0.93%  "
$ws.Range("D26").Value = "0.1414"
$ws.Range("E26").Value = "  -1.07%  "
$ws.Range("D27").Value = "8.525"
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("D28").Value = "17.85"
$ws.Range("E28").Value = "  -0.56%  "
$ws.Range("D29").Value = "1.495"
$ws.Range("E29").Value = "  -0.35%  "
$ws.Range("D30").Value = "0.05984"
$ws.Range("E30").Value = "  -3.83%  "
$ws.Range("D31").Value = "1.257"
$ws.Range("E31").Value = "  -1.28%  "
$ws.Range("D32").Value = "4.139"
$ws.Range("E32").Value = "  -0.22%  "
$ws.Range("D33").Value = "4.087"
$ws.Range("E33").Value = "  -0.75%  "
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("D35").Value = "1.177"
$ws.Range("E35").Value = "  +1.38%  "
$ws.Range("D36").Value = "0.7239"
$ws.Range("E36").Value = "  -1.16%  "
$ws.Range("D37").Value = "2.602"
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("D38").Value = "2.783"
$ws.Range("E38").Value = "  -2.60%  "
$ws.Range("D39").Value = "0.01780"
$ws.Range("E39").Value = "  -1.08%  "
$ws.Range("D40").Value = "1.207.77"
$ws.Range("E40").Value = "  -0.48%  "
$ws.Range("D41").Value = "0.9127"
$ws.Range("E41").Value = "  -1.41%  "
$ws.Range("D42").Value = "6.195"
$ws.Range("E42").Value = "  +0.35%  "
$ws.Range("D43").Value = "2.075.87"
$ws.Range("E43").Value = "  +2.82%  "
$ws.Range("D44").Value = "0.9992"
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("D45").Value = "102.04"
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").Value = "67.18"
$ws.Range("E46").Value = "  +0.99%  "
$ws.Range("D47").Value = "7.373"
$ws.Range("E47").Value = "  +10.02%  "
$ws.Range("D50").Value = "9.153"
$ws.Range("E50").Value = "  -1.75%  "
$ws.Range("D51").Value = "1.667"
$ws.Range("E51").Value = "  +1.90%  "

# Row 22/23 swap: Dai <-> Chainlink (ranking reorder)
# Row 48/49 swap: BabyDogeCoin <-> TheSandbox (ranking reorder)
$ws.Range("B22").Value = "Chainlink"
$ws.Range("C22").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D22").Value = "7.768"
$ws.Range("E22").Value = "  +4.59%  "
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").Value = "0.9993"
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("B48").Value = "TheSandbox"
$ws.Range("C48").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D48").Value = "0.4064"
$ws.Range("E48").Value = "  -0.51%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.00000000118"
$ws.Range("E49").Value = "  -1.71%  "

# Restore the default (un-styled) cell appearance now that the values are
# safely stored as text, matching the original workbook styling.
$ws.Range("D2:D51").Style = "Normal"
